# Convert to excel from multiple arb file
#
# This localization sheet previously held only an English ("en") value
# column. Re-generating it from multiple .arb files adds a Polish ("pl")
# value column, inserts a new "green_color" entry, and appends a new
# "pink_color" entry (which only has a Polish translation).
#
# Cells are (re)written row by row, column by column, matching the order
# in which a real generator would emit them from the source .arb files.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - headers
$ws.Range("A1").Value = "category"
$ws.Range("B1").Value = "text"
$ws.Range("C1").Value = "description"
$ws.Range("D1").Value = "en"
$ws.Range("E1").Value = "pl"

# Row 2 - @@locale
$ws.Range("B2").Value = "@@locale"
$ws.Range("D2").Value = "en"
$ws.Range("E2").Value = "pl"

# Row 3 - company_name
$ws.Range("B3").Value = "company_name"
$ws.Range("D3").Value = "LLC"
$ws.Range("E3").Value = "S.A."

# Row 4 - red_color
$ws.Range("B4").Value = "red_color"
$ws.Range("D4").Value = "Red"
$ws.Range("E4").Value = "Czerwony"

# Row 5 - green_color (new entry, English only initially)
$ws.Range("B5").Value = "green_color"
$ws.Range("D5").Value = "Green"

# Row 6 - orange_color
$ws.Range("B6").Value = "orange_color"
$ws.Range("D6").Value = "Orange"
$ws.Range("E6").Value = "Pomarańczowy"

# Row 7 - pink_color (new entry, Polish only)
$ws.Range("B7").Value = "pink_color"
$ws.Range("E7").Value = "Różowy"

# Give the new "pl" column (E) the same auto-fit style width as column B
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(2).ColumnWidth
